$wb = $excel.ActiveWorkbook

# The dataset_internal_id and indicator_internal_id values live on the
# "metadata" sheet, in column B, next to their respective labels in column A.
$ws = $wb.Worksheets.Item("metadata")

# B2 holds dataset_internal_id ("LandAndGender" -> "LG")
$ws.Range("B2").Value = "LG"

# B3 holds indicator_internal_id ("LandAndGender.3Fa" -> "LG.3Fa")
$ws.Range("B3").Value = "LG.3Fa"
